# Apply the "feat: add 2022-Q4 data" edit:
#  1. Insert a new row 2 into the "总计" (totals) sheet with the 2022-Q4
#     summary figures, shifting the existing quarters down by one row.
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计",
#     holding the per-fund holdings detail for the new quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# --- 1. "总计" sheet: insert new row 2 for 2022-Q4 -------------------------
$total.Rows.Item(2).Insert()
# The inserted row inherits formatting from the row above (the bold header);
# strip that back to the plain (unstyled) look used by the other data rows.
$total.Range("A2:D2").ClearFormats()
# Re-apply the same label style ("s=2") used by the other index cells in
# column A, by copying format from the cell directly below.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 25
$total.Range("D2").Value = 2.8

# --- 2. New "2022-Q4" worksheet with fund holdings detail ------------------
$new = $wb.Worksheets.Add($null, $total)
$new.Name = "2022-Q4"

# Give the header row (B1:H1) and the index column (A2:A26) the bold/bordered
# "label" style used throughout the workbook, copied from the "总计" sheet so
# we reuse the existing style index instead of fabricating a new one.
$total.Range("B1:D1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2:A8").Copy()
$new.Range("A2:A26").PasteSpecial(-4122)

$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"
$new.Range("A2").Value = 0
$new.Range("B2").Value = "'010699"
$new.Range("C2").Value = "东方红创新趋势混合"
$new.Range("D2").Value = "'24.58"
$new.Range("E2").Value = "'89.62"
$new.Range("F2").Value = "'3.23"
$new.Range("G2").Value = "'0.7939"
$new.Range("H2").Value = 8
$new.Range("A3").Value = 1
$new.Range("B3").Value = "'008271"
$new.Range("C3").Value = "大成优势企业混合A"
$new.Range("D3").Value = "'9.34"
$new.Range("E3").Value = "'79.90"
$new.Range("F3").Value = "'6.45"
$new.Range("G3").Value = "'0.6024"
$new.Range("H3").Value = 5
$new.Range("A4").Value = 2
$new.Range("B4").Value = "'001487"
$new.Range("C4").Value = "宝盈优势产业灵活配置混合A"
$new.Range("D4").Value = "'10.62"
$new.Range("E4").Value = "'94.48"
$new.Range("F4").Value = "'3.14"
$new.Range("G4").Value = "'0.3335"
$new.Range("H4").Value = 7
$new.Range("A5").Value = 3
$new.Range("B5").Value = "'001128"
$new.Range("C5").Value = "宝盈新兴产业灵活配置混合A"
$new.Range("D5").Value = "'9.13"
$new.Range("E5").Value = "'94.44"
$new.Range("F5").Value = "'3.05"
$new.Range("G5").Value = "'0.2785"
$new.Range("H5").Value = 9
$new.Range("A6").Value = 4
$new.Range("B6").Value = "'001877"
$new.Range("C6").Value = "宝盈国家安全沪港深股票A"
$new.Range("D6").Value = "'6.58"
$new.Range("E6").Value = "'94.32"
$new.Range("F6").Value = "'3.87"
$new.Range("G6").Value = "'0.2546"
$new.Range("H6").Value = 9
$new.Range("A7").Value = 5
$new.Range("B7").Value = "'009069"
$new.Range("C7").Value = "大成睿鑫股票A"
$new.Range("D7").Value = "'3.39"
$new.Range("E7").Value = "'82.20"
$new.Range("F7").Value = "'5.04"
$new.Range("G7").Value = "'0.1709"
$new.Range("H7").Value = 8
$new.Range("A8").Value = 6
$new.Range("B8").Value = "'012771"
$new.Range("C8").Value = "宝盈优势产业灵活配置混合C"
$new.Range("D8").Value = "'3.19"
$new.Range("E8").Value = "'94.48"
$new.Range("F8").Value = "'3.14"
$new.Range("G8").Value = "'0.1002"
$new.Range("H8").Value = 7
$new.Range("A9").Value = 7
$new.Range("B9").Value = "'002707"
$new.Range("C9").Value = "摩根士丹利华鑫科技领先灵活配置混合A"
$new.Range("D9").Value = "'1.76"
$new.Range("E9").Value = "'92.23"
$new.Range("F9").Value = "'3.37"
$new.Range("G9").Value = "'0.0593"
$new.Range("H9").Value = 7
$new.Range("A10").Value = 8
$new.Range("B10").Value = "'008272"
$new.Range("C10").Value = "大成优势企业混合C"
$new.Range("D10").Value = "'0.87"
$new.Range("E10").Value = "'79.90"
$new.Range("F10").Value = "'6.45"
$new.Range("G10").Value = "'0.0561"
$new.Range("H10").Value = 5
$new.Range("A11").Value = 9
$new.Range("B11").Value = "'012815"
$new.Range("C11").Value = "宝盈新兴产业灵活配置混合C"
$new.Range("D11").Value = "'1.06"
$new.Range("E11").Value = "'94.44"
$new.Range("F11").Value = "'3.05"
$new.Range("G11").Value = "'0.0323"
$new.Range("H11").Value = 9
$new.Range("A12").Value = 10
$new.Range("B12").Value = "'002103"
$new.Range("C12").Value = "招商康泰灵活配置混合"
$new.Range("D12").Value = "'1.43"
$new.Range("E12").Value = "'39.33"
$new.Range("F12").Value = "'1.61"
$new.Range("G12").Value = "'0.0230"
$new.Range("H12").Value = 9
$new.Range("A13").Value = 11
$new.Range("B13").Value = "'009070"
$new.Range("C13").Value = "大成睿鑫股票C"
$new.Range("D13").Value = "'0.36"
$new.Range("E13").Value = "'82.20"
$new.Range("F13").Value = "'5.04"
$new.Range("G13").Value = "'0.0181"
$new.Range("H13").Value = 8
$new.Range("A14").Value = 12
$new.Range("B14").Value = "'006573"
$new.Range("C14").Value = "人保行业轮动混合A"
$new.Range("D14").Value = "'0.71"
$new.Range("E14").Value = "'79.80"
$new.Range("F14").Value = "'2.34"
$new.Range("G14").Value = "'0.0166"
$new.Range("H14").Value = 5
$new.Range("A15").Value = 13
$new.Range("B15").Value = "'570007"
$new.Range("C15").Value = "诺德优选30混合"
$new.Range("D15").Value = "'0.19"
$new.Range("E15").Value = "'87.62"
$new.Range("F15").Value = "'7.65"
$new.Range("G15").Value = "'0.0145"
$new.Range("H15").Value = 2
$new.Range("A16").Value = 14
$new.Range("B16").Value = "'008082"
$new.Range("C16").Value = "国寿安保研究精选混合A"
$new.Range("D16").Value = "'0.34"
$new.Range("E16").Value = "'84.56"
$new.Range("F16").Value = "'3.06"
$new.Range("G16").Value = "'0.0104"
$new.Range("H16").Value = 9
$new.Range("A17").Value = 15
$new.Range("B17").Value = "'010765"
$new.Range("C17").Value = "国寿安保华丰混合A"
$new.Range("D17").Value = "'0.40"
$new.Range("E17").Value = "'83.84"
$new.Range("F17").Value = "'2.27"
$new.Range("G17").Value = "'0.0091"
$new.Range("H17").Value = 6
$new.Range("A18").Value = 16
$new.Range("B18").Value = "'007316"
$new.Range("C18").Value = "交银施罗德可转债债券A"
$new.Range("D18").Value = "'0.82"
$new.Range("E18").Value = "'20.54"
$new.Range("F18").Value = "'0.86"
$new.Range("G18").Value = "'0.0071"
$new.Range("H18").Value = 2
$new.Range("A19").Value = 17
$new.Range("B19").Value = "'013613"
$new.Range("C19").Value = "宝盈国家安全沪港深股票C"
$new.Range("D19").Value = "'0.13"
$new.Range("E19").Value = "'94.32"
$new.Range("F19").Value = "'3.87"
$new.Range("G19").Value = "'0.0050"
$new.Range("H19").Value = 9
$new.Range("A20").Value = 18
$new.Range("B20").Value = "'010762"
$new.Range("C20").Value = "博时恒康一年持有期混合A"
$new.Range("D20").Value = "'0.45"
$new.Range("E20").Value = "'21.62"
$new.Range("F20").Value = "'0.95"
$new.Range("G20").Value = "'0.0043"
$new.Range("H20").Value = 8
$new.Range("A21").Value = 19
$new.Range("B21").Value = "'008083"
$new.Range("C21").Value = "国寿安保研究精选混合C"
$new.Range("D21").Value = "'0.13"
$new.Range("E21").Value = "'84.56"
$new.Range("F21").Value = "'3.06"
$new.Range("G21").Value = "'0.0040"
$new.Range("H21").Value = 9
$new.Range("A22").Value = 20
$new.Range("B22").Value = "'006574"
$new.Range("C22").Value = "人保行业轮动混合C"
$new.Range("D22").Value = "'0.16"
$new.Range("E22").Value = "'79.80"
$new.Range("F22").Value = "'2.34"
$new.Range("G22").Value = "'0.0037"
$new.Range("H22").Value = 5
$new.Range("A23").Value = 21
$new.Range("B23").Value = "'007317"
$new.Range("C23").Value = "交银施罗德可转债债券C"
$new.Range("D23").Value = "'0.22"
$new.Range("E23").Value = "'20.54"
$new.Range("F23").Value = "'0.86"
$new.Range("G23").Value = "'0.0019"
$new.Range("H23").Value = 2
$new.Range("A24").Value = 22
$new.Range("B24").Value = "'014871"
$new.Range("C24").Value = "摩根士丹利华鑫科技领先灵活配置混合C"
$new.Range("D24").Value = "'0.05"
$new.Range("E24").Value = "'92.23"
$new.Range("F24").Value = "'3.37"
$new.Range("G24").Value = "'0.0017"
$new.Range("H24").Value = 7
$new.Range("A25").Value = 23
$new.Range("B25").Value = "'010763"
$new.Range("C25").Value = "博时恒康一年持有期混合C"
$new.Range("D25").Value = "'0.11"
$new.Range("E25").Value = "'21.62"
$new.Range("F25").Value = "'0.95"
$new.Range("G25").Value = "'0.0010"
$new.Range("H25").Value = 8
$new.Range("A26").Value = 24
$new.Range("B26").Value = "'010766"
$new.Range("C26").Value = "国寿安保华丰混合C"
$new.Range("D26").Value = "'0.01"
$new.Range("E26").Value = "'83.84"
$new.Range("F26").Value = "'2.27"
$new.Range("G26").Value = "'0.0002"
$new.Range("H26").Value = 6

# Several of the numeric-looking columns (fund code, fund size, position
# weight, ...) are stored as TEXT in this workbook (leading zeros in fund
# codes, fixed trailing zeros in percentages, etc.). Entering them with a
# leading apostrophe forces Excel to keep them as text instead of silently
# converting to numbers; that also flags the cells with a "quotePrefix"
# style, so finish by pasting the (unstyled) format from a blank cell over
# just those text columns to land on the same plain style as the rest of
# the sheet.
$new.Range("Z100").Copy()
$new.Range("B2:B26").PasteSpecial(-4122)
$new.Range("D2:G26").PasteSpecial(-4122)
